$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 375.14285
$ws.Range("I4").Value = 270.83334
$ws.Range("K4").Value = 270.83334
$ws.Range("M4").Value = -156.83334
$ws.Range("H18").Value = 1596.2858
$ws.Range("J18").Value = 699.5
$ws.Range("L18").Value = 699.5
$ws.Range("N18").Value = -1267.5
$ws.Range("H39").Value = 46.545456
$ws.Range("I39").Value = 46.545456
$ws.Range("K39").Value = 139.636368
$ws.Range("M39").Value = 156.363632
$ws.Range("H80").Value = 530
$ws.Range("I80").Value = 392.85715
$ws.Range("K80").Value = 1178.57145
$ws.Range("M80").Value = -180.5714499999999
$ws.Range("H83").Value = 530
$ws.Range("I83").Value = 392.85715
$ws.Range("K83").Value = 3535.71435
$ws.Range("M83").Value = 1456.28565
$ws.Range("H86").Value = 5069.9
$ws.Range("I86").Value = 3066.3333
$ws.Range("K86").Value = 3066.3333
$ws.Range("M86").Value = -1943.3333
$ws.Range("H89").Value = 5069.9
$ws.Range("I89").Value = 3066.3333
$ws.Range("K89").Value = 15331.6665
$ws.Range("M89").Value = -9715.666499999999
$ws.Range("H111").Value = 3717.4285
$ws.Range("I111").Value = 4005.4
$ws.Range("K111").Value = 12016.2
$ws.Range("M111").Value = -8949.200000000001
$ws.Range("H125").Value = 988.2
$ws.Range("I125").Value = 986
$ws.Range("J125").Value = 991.5
$ws.Range("K125").Value = 8874
$ws.Range("L125").Value = 8923.5
$ws.Range("M125").Value = -6414
$ws.Range("N125").Value = -13843.5
$ws.Range("H138").Value = 4194.1333
$ws.Range("J138").Value = 4538.8716
$ws.Range("L138").Value = 13616.6148
$ws.Range("N138").Value = -23896.6148

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2499
$ws.Range("I61").Value = 2499
$ws.Range("K61").Value = 2499
$ws.Range("M61").Value = -2287
$ws.Range("H63").Value = 7184.5264
$ws.Range("I63").Value = 7900
$ws.Range("J63").Value = 6929
$ws.Range("K63").Value = 7900
$ws.Range("L63").Value = 6929
$ws.Range("M63").Value = -7214
$ws.Range("N63").Value = -8301
$ws.Range("H66").Value = 7184.5264
$ws.Range("I66").Value = 7900
$ws.Range("J66").Value = 6929
$ws.Range("K66").Value = 39500
$ws.Range("L66").Value = 34645
$ws.Range("M66").Value = -36068
$ws.Range("N66").Value = -41509
$ws.Range("H132").Value = 2019.5264
$ws.Range("I132").Value = 2021.9412
$ws.Range("K132").Value = 6065.8236
$ws.Range("M132").Value = -3535.8236
$ws.Range("H133").Value = 67473.25
$ws.Range("J133").Value = 67473.25
$ws.Range("L133").Value = 67473.25
$ws.Range("N133").Value = -72533.25
$ws.Range("H136").Value = 2499
$ws.Range("I136").Value = 2499
$ws.Range("K136").Value = 7497
$ws.Range("M136").Value = -4947

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1520.25
$ws.Range("I86").Value = 1575
$ws.Range("K86").Value = 1575
$ws.Range("M86").Value = -452
$ws.Range("H89").Value = 1520.25
$ws.Range("I89").Value = 1575
$ws.Range("K89").Value = 7875
$ws.Range("M89").Value = -2259
$ws.Range("H99").Value = 3533.2173
$ws.Range("I99").Value = 3452.7646
$ws.Range("K99").Value = 3452.7646
$ws.Range("M99").Value = -1954.7646
$ws.Range("H134").Value = 1597.9412
$ws.Range("I134").Value = 1588.4375
$ws.Range("K134").Value = 4765.3125
$ws.Range("M134").Value = -2230.3125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 9598.777
$ws.Range("I86").Value = 8678
$ws.Range("K86").Value = 8678
$ws.Range("M86").Value = -7555
$ws.Range("H89").Value = 9598.777
$ws.Range("I89").Value = 8678
$ws.Range("K89").Value = 43390
$ws.Range("M89").Value = -37774
$ws.Range("H132").Value = 2606.125
$ws.Range("I132").Value = 1485.1
$ws.Range("K132").Value = 4455.299999999999
$ws.Range("M132").Value = -1925.299999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55645.445
$ws.Range("I2").Value = 71499.71000000001
$ws.Range("K2").Value = 428998.26
$ws.Range("M2").Value = -428885.26
$ws.Range("H62").Value = 1959.8
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 6000
$ws.Range("H65").Value = 1959.8
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 18000
$ws.Range("H92").Value = 1000
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("H140").Value = 2555.9092
$ws.Range("I140").Value = 3064.375
$ws.Range("K140").Value = 9193.125
$ws.Range("M140").Value = -4013.125
$ws.Range("H141").Value = 5000
$ws.Range("I141").Value = 5000
$ws.Range("K141").Value = 15000
$ws.Range("N62").Value = -7372
$ws.Range("N65").Value = -24864
$ws.Range("N92").Value = -5496
$ws.Range("M141").Value = -9820

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5499.5
$ws.Range("I70").Value = 6499
$ws.Range("K70").Value = 6499
$ws.Range("M70").Value = -6229
$ws.Range("H73").Value = 5499.5
$ws.Range("I73").Value = 6499
$ws.Range("K73").Value = 6499
$ws.Range("M73").Value = -5563
$ws.Range("H107").Value = 1409.9445
$ws.Range("J107").Value = 968
$ws.Range("L107").Value = 968
$ws.Range("N107").Value = -4808
$ws.Range("H126").Value = 4213
$ws.Range("I126").Value = 3247.1667
$ws.Range("K126").Value = 9741.500100000001
$ws.Range("M126").Value = -7271.500100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3648.8845
$ws.Range("I132").Value = 2713.1177
$ws.Range("K132").Value = 8139.353099999999
$ws.Range("M132").Value = -5609.353099999999
$ws.Range("H136").Value = 2637.5
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 2683.3333
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 8049.999899999999
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -13149.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8199.333000000001
$ws.Range("J62").Value = 8537.77
$ws.Range("L62").Value = 8537.77
$ws.Range("N62").Value = -9785.77
$ws.Range("H65").Value = 8199.333000000001
$ws.Range("J65").Value = 8537.77
$ws.Range("L65").Value = 42688.85000000001
$ws.Range("N65").Value = -48928.85000000001
$ws.Range("H132").Value = 2034.625
$ws.Range("I132").Value = 1350.3846
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 4051.1538
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -1521.1538
$ws.Range("N132").Value = -20058.9995
